# Updates cryptos list: prices (D) and 1h volume % (E) refreshed;
# rows 34/35 (NEARProtocol <-> Kaspa) and 50/51 (ONDO <-> Cosmos) swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.482.98"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "'3.675.31"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'640.49"
$ws.Range("E5").Value = "  -5.71%  "

$ws.Range("D6").Value = "'159.43"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.497"
$ws.Range("E8").Value = "  +0.56%  "

$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").Value = "'7.10"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("E11").Value = "  +1.24%  "

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "'4.296.49"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "'32.69"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").Value = "'3.687.39"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "'69.476.78"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "'15.97"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").Value = "'6.49"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").Value = "'466.78"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").Value = "'79.29"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").Value = "'3.823.26"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  +1.67%  "

$ws.Range("D27").Value = "'10.88"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("D28").Value = "'9.03"
$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("D29").Value = "'2.62"
$ws.Range("E29").Value = "  -2.72%  "

$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "'26.84"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.165"
$ws.Range("E34").Value = "  +3.56%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.46"
$ws.Range("E35").Value = "  -1.48%  "

$ws.Range("D36").Value = "'3.667.44"
$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").Value = "'8.46"
$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'5.88"
$ws.Range("E39").Value = "  -6.87%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").Value = "'177.13"
$ws.Range("E41").Value = "  +4.40%  "

$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("D43").Value = "'0.0899"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").Value = "'0.925"
$ws.Range("E44").Value = "  -1.73%  "

$ws.Range("D45").Value = "'46.83"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("D47").Value = "'27.42"
$ws.Range("E47").Value = "  -1.95%  "

$ws.Range("D48").Value = "'0.000269"
$ws.Range("E48").Value = "  -4.57%  "

$ws.Range("D49").Value = "'1.07"
$ws.Range("E49").Value = "  -3.36%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.84"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'1.25"
$ws.Range("E51").Value = "  -3.17%  "
